$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 7.038349
$ws.Range("H2").Value = 21.115047
$ws.Range("I2").Value = 0.8947789352175559
$ws.Range("J2").Value = 0.894778935217556
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.248835333333334
$ws.Range("N2").Value = 6.746506
$ws.Range("O2").Value = 0.03590294220158827
$ws.Range("P2").Value = 0.03590294220158827
$ws.Range("Q2").Value = 15.82808791953133
$ws.Range("R2").Value = 142.452791275782
$ws.Range("S2").Value = 0.03212519639431461
$ws.Range("T2").Value = 0.03212519639431461

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 7.038349
$ws.Range("H3").Value = 21.115047
$ws.Range("I3").Value = 0.8947789352175559
$ws.Range("J3").Value = 0.894778935217556
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 44.29005966666667
$ws.Range("N3").Value = 132.870179
$ws.Range("O3").Value = 0.7070964373190639
$ws.Range("P3").Value = 0.7070964373190639
$ws.Range("Q3").Value = 311.7288971648237
$ws.Range("R3").Value = 2805.560074483413
$ws.Range("S3").Value = 0.6326949972804793
$ws.Range("T3").Value = 0.6326949972804794

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 7.038349
$ws.Range("H4").Value = 21.115047
$ws.Range("I4").Value = 0.8947789352175559
$ws.Range("J4").Value = 0.894778935217556
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 16.09762433333333
$ws.Range("N4").Value = 48.292873
$ws.Range("O4").Value = 0.2570006204793478
$ws.Range("P4").Value = 0.2570006204793479
$ws.Range("Q4").Value = 113.3006981288923
$ws.Range("R4").Value = 1019.706283160031
$ws.Range("S4").Value = 0.229958741542762
$ws.Range("T4").Value = 0.2299587415427621

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.09030199999999999
$ws.Range("H5").Value = 0.270906
$ws.Range("I5").Value = 0.01148001149247014
$ws.Range("J5").Value = 0.01148001149247014
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.248835333333334
$ws.Range("N5").Value = 6.746506
$ws.Range("O5").Value = 0.03590294220158827
$ws.Range("P5").Value = 0.03590294220158827
$ws.Range("Q5").Value = 0.2030743282706667
$ws.Range("R5").Value = 1.827668954436
$ws.Range("S5").Value = 0.0004121661890877246
$ws.Range("T5").Value = 0.0004121661890877247

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.09030199999999999
$ws.Range("H6").Value = 0.270906
$ws.Range("I6").Value = 0.01148001149247014
$ws.Range("J6").Value = 0.01148001149247014
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 44.29005966666667
$ws.Range("N6").Value = 132.870179
$ws.Range("O6").Value = 0.7070964373190639
$ws.Range("P6").Value = 0.7070964373190639
$ws.Range("Q6").Value = 3.999480968019333
$ws.Range("R6").Value = 35.995328712174
$ws.Range("S6").Value = 0.008117475226707546
$ws.Range("T6").Value = 0.008117475226707548

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.09030199999999999
$ws.Range("H7").Value = 0.270906
$ws.Range("I7").Value = 0.01148001149247014
$ws.Range("J7").Value = 0.01148001149247014
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 16.09762433333333
$ws.Range("N7").Value = 48.292873
$ws.Range("O7").Value = 0.2570006204793478
$ws.Range("P7").Value = 0.2570006204793479
$ws.Range("Q7").Value = 1.453647672548666
$ws.Range("R7").Value = 13.082829052938
$ws.Range("S7").Value = 0.00295037007667487
$ws.Range("T7").Value = 0.002950370076674871

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.4911496666666667
$ws.Range("H8").Value = 1.473449
$ws.Range("I8").Value = 0.06243941239237462
$ws.Range("J8").Value = 0.06243941239237463
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.248835333333334
$ws.Range("N8").Value = 6.746506
$ws.Range("O8").Value = 0.03590294220158827
$ws.Range("P8").Value = 0.03590294220158827
$ws.Range("Q8").Value = 1.104514724354889
$ws.Range("R8").Value = 9.940632519194001
$ws.Range("S8").Value = 0.00224175861422456
$ws.Range("T8").Value = 0.002241758614224561

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.4911496666666667
$ws.Range("H9").Value = 1.473449
$ws.Range("I9").Value = 0.06243941239237462
$ws.Range("J9").Value = 0.06243941239237463
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 44.29005966666667
$ws.Range("N9").Value = 132.870179
$ws.Range("O9").Value = 0.7070964373190639
$ws.Range("P9").Value = 0.7070964373190639
$ws.Range("Q9").Value = 21.75304804193011
$ws.Range("R9").Value = 195.777432377371
$ws.Range("S9").Value = 0.0441506860509439
$ws.Range("T9").Value = 0.0441506860509439

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.4911496666666667
$ws.Range("H10").Value = 1.473449
$ws.Range("I10").Value = 0.06243941239237462
$ws.Range("J10").Value = 0.06243941239237463
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 16.09762433333333
$ws.Range("N10").Value = 48.292873
$ws.Range("O10").Value = 0.2570006204793478
$ws.Range("P10").Value = 0.2570006204793479
$ws.Range("Q10").Value = 7.906342825441888
$ws.Range("R10").Value = 71.157085428977
$ws.Range("S10").Value = 0.01604696772720616
$ws.Range("T10").Value = 0.01604696772720616

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.1894393333333333
$ws.Range("H11").Value = 0.568318
$ws.Range("I11").Value = 0.02408325091130372
$ws.Range("J11").Value = 0.02408325091130372
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 2.248835333333334
$ws.Range("N11").Value = 6.746506
$ws.Range("O11").Value = 0.03590294220158827
$ws.Range("P11").Value = 0.03590294220158827
$ws.Range("Q11").Value = 0.4260178663231111
$ws.Range("R11").Value = 3.834160796908
$ws.Range("S11").Value = 0.0008646595654948856
$ws.Range("T11").Value = 0.0008646595654948858

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.1894393333333333
$ws.Range("H12").Value = 0.568318
$ws.Range("I12").Value = 0.02408325091130372
$ws.Range("J12").Value = 0.02408325091130372
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 44.29005966666667
$ws.Range("N12").Value = 132.870179
$ws.Range("O12").Value = 0.7070964373190639
$ws.Range("P12").Value = 0.7070964373190639
$ws.Range("Q12").Value = 8.390279376546889
$ws.Range("R12").Value = 75.512514388922
$ws.Range("S12").Value = 0.01702918091844396
$ws.Range("T12").Value = 0.01702918091844396

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.1894393333333333
$ws.Range("H13").Value = 0.568318
$ws.Range("I13").Value = 0.02408325091130372
$ws.Range("J13").Value = 0.02408325091130372
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 16.09762433333333
$ws.Range("N13").Value = 48.292873
$ws.Range("O13").Value = 0.2570006204793478
$ws.Range("P13").Value = 0.2570006204793479
$ws.Range("Q13").Value = 3.049523221957111
$ws.Range("R13").Value = 27.445708997614
$ws.Range("S13").Value = 0.006189410427364875
$ws.Range("T13").Value = 0.006189410427364876

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.05678
$ws.Range("H14").Value = 0.17034
$ws.Range("I14").Value = 0.007218389986295482
$ws.Range("J14").Value = 0.007218389986295483
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 2.248835333333334
$ws.Range("N14").Value = 6.746506
$ws.Range("O14").Value = 0.03590294220158827
$ws.Range("P14").Value = 0.03590294220158827
$ws.Range("Q14").Value = 0.1276888702266667
$ws.Range("R14").Value = 1.14919983204
$ws.Range("S14").Value = 0.0002591614384664903
$ws.Range("T14").Value = 0.0002591614384664903

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.05678
$ws.Range("H15").Value = 0.17034
$ws.Range("I15").Value = 0.007218389986295482
$ws.Range("J15").Value = 0.007218389986295483
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 44.29005966666667
$ws.Range("N15").Value = 132.870179
$ws.Range("O15").Value = 0.7070964373190639
$ws.Range("P15").Value = 0.7070964373190639
$ws.Range("Q15").Value = 2.514789587873334
$ws.Range("R15").Value = 22.63310629086
$ws.Range("S15").Value = 0.005104097842489142
$ws.Range("T15").Value = 0.005104097842489143

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.05678
$ws.Range("H16").Value = 0.17034
$ws.Range("I16").Value = 0.007218389986295482
$ws.Range("J16").Value = 0.007218389986295483
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 16.09762433333333
$ws.Range("N16").Value = 48.292873
$ws.Range("O16").Value = 0.2570006204793478
$ws.Range("P16").Value = 0.2570006204793479
$ws.Range("Q16").Value = 0.9140231096466666
$ws.Range("R16").Value = 8.22620798682
$ws.Range("S16").Value = 0.00185513070533985
$ws.Range("T16").Value = 0.001855130705339851
